$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add a new row of data (row 25) below the existing table
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = 3
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = -1000
$ws.Range("H25").Value = 3000
$ws.Range("I25").Value = -2000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 98
$ws.Range("M25").Value = "Outside"
$ws.Range("N25").Value = 0
$ws.Range("O25").Value = 4140

# Update the current selection to match the new cursor position
$ws.Range("I26").Select()
